$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 18:16"

# --- Countries whose statistics were refreshed (values only, A column untouched) ---
$ws.Range("B4").Value = 3180200
$ws.Range("C4").Value = 21268
$ws.Range("D4").Value = 1394249
$ws.Range("E4").Value = 1650740
$ws.Range("G4").Value = 349
$ws.Range("H4").Value = 135211

$ws.Range("B6").Value = 791001
$ws.Range("C6").Value = 21949
$ws.Range("D6").Value = 493484
$ws.Range("E6").Value = 275925
$ws.Range("G6").Value = 448
$ws.Range("H6").Value = 21592

$ws.Range("B9").Value = 306216
$ws.Range("C9").Value = 3133
$ws.Range("D9").Value = 274922
$ws.Range("E9").Value = 24612
$ws.Range("G9").Value = 109
$ws.Range("H9").Value = 6682

$ws.Range("B14").Value = 242363
$ws.Range("C14").Value = 214
$ws.Range("D14").Value = 193978
$ws.Range("E14").Value = 13459
$ws.Range("G14").Value = 12
$ws.Range("H14").Value = 34926

$ws.Range("B23").Value = 106741
$ws.Range("C23").Value = 307
$ws.Range("D23").Value = 70503
$ws.Range("E23").Value = 27492
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 8746

$ws.Range("B28").Value = 74333
$ws.Range("C28").Value = 48
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = 5500

$ws.Range("D41").Value = 41645
$ws.Range("E41").Value = 3752

$ws.Range("B58").Value = 22464
$ws.Range("C58").Value = 548
$ws.Range("D58").Value = 13591
$ws.Range("E58").Value = 8589
$ws.Range("G58").Value = 10
$ws.Range("H58").Value = 284

$ws.Range("B69").Value = 12859
$ws.Range("C69").Value = 45
$ws.Range("D69").Value = 8123
$ws.Range("E69").Value = 4384
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 352

$ws.Range("B86").Value = 6410
$ws.Range("C86").Value = 46
$ws.Range("D86").Value = 5067
$ws.Range("E86").Value = 1289

$ws.Range("B100").Value = 3672
$ws.Range("C100").Value = 50
$ws.Range("E100").Value = 2105

$ws.Range("D133").Value = 982
$ws.Range("E133").Value = 177

# --- Montenegro jumps above "Republica de Chipre" in the ranking; rows 138-143
#     shift down by one and row 138 receives Montenegro's refreshed data ---
$ws.Range("A138").Value = "Montenegro"
$ws.Range("B138").Value = 1019
$ws.Range("C138").Value = 59
$ws.Range("D138").Value = 320
$ws.Range("E138").Value = 680
$ws.Range("G138").Value = 2
$ws.Range("H138").Value = 19

$ws.Range("A139").Value = "Republica de Chipre"
$ws.Range("B139").Value = 1008
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 839
$ws.Range("E139").Value = 150
$ws.Range("H139").Value = 19

$ws.Range("A140").Value = "Burkina Faso"
$ws.Range("B140").Value = 1005
$ws.Range("C140").Value = 2
$ws.Range("D140").Value = 862
$ws.Range("E140").Value = 90
$ws.Range("H140").Value = 53

$ws.Range("A141").Value = "Uganda"
$ws.Range("B141").Value = 1000
$ws.Range("C141").Value = 23
$ws.Range("D141").Value = 908
$ws.Range("E141").Value = 92
$ws.Range("H141").Value = 0

$ws.Range("A142").Value = "Uruguay"
$ws.Range("B142").Value = 974
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 871
$ws.Range("E142").Value = 74
$ws.Range("H142").Value = 29

$ws.Range("A143").Value = "Georgia"
$ws.Range("B143").Value = 968
$ws.Range("C143").Value = 5
$ws.Range("D143").Value = 844
$ws.Range("E143").Value = 109
$ws.Range("H143").Value = 15

# --- Tied countries re-sorted alphabetically; data identical, only names swap ---
$ws.Range("A184").Value = "Seychelles"
$ws.Range("A185").Value = "Lesoto"

$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
